# Generate Report for Handoff
# Updates the "b.md" row across all three sheets to reflect a fresh handoff
# (new xliff files generated, new timestamps, and a new "stale handback"
# error detail raised for b.md), and fixes the de-de "a.md" row's handoff
# file name which was erroneously suffixed with ".zh-cn.xlf" instead of
# ".de-de.xlf".

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9423d9c9e9f522a22a1bf4e5011c45de41e37a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df2091be9c5658a141f6fbed4109fc53bbfbbd48/e2e/b.md."

# ---- Overview sheet: b.md row (row 3) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "'Ready for handoff"
$wsOverview.Range("F3").Value = "'Ready for handoff"
$wsOverview.Range("G3").Value = "'2016-09-03 16:44:17"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# b.md row (row 3): new handoff generated
$wsZhCn.Range("C3").Value = "'Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "'b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "'2016-09-03 16:44:13"
$wsZhCn.Range("P3").Value = "'" + $errorDetail

# Widen the Error Detail column (P) so the new message is readable
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

# a.md row (row 2): fix mislabeled handoff file language suffix + handback datetime
$wsDeDe.Range("G2").Value = "'a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "'2016-09-03 16:43:54"

# b.md row (row 3): new handoff generated
$wsDeDe.Range("C3").Value = "'Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "'b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "'2016-09-03 16:44:17"
$wsDeDe.Range("J3").Value = "'a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "'2016-09-03 16:43:54"
$wsDeDe.Range("P3").Value = "'" + $errorDetail

# Widen the Error Detail column (P) so the new message is readable
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664

Write-Output "Generate Report for Handoff: applied"
